$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.926.95'
$ws.Range('E2').Value = '  +2.69%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.598.53'
$ws.Range('E3').Value = '  +2.72%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.09'
$ws.Range('E5').Value = '  +2.61%  '
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('E7').Value = '  +1.11%  '
$ws.Range('E8').Value = '  +1.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0612'
$ws.Range('E9').Value = '  +0.57%  '
$ws.Range('E10').Value = '  +1.05%  '
$ws.Range('E11').Value = '  +3.90%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.822.05'
$ws.Range('E12').Value = '  +2.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.597.08'
$ws.Range('E13').Value = '  +2.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.508'
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.958.64'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.17'
$ws.Range('E17').Value = '  +2.00%  '
$ws.Range('E18').Value = '  +1.91%  '
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '200.44'
$ws.Range('E20').Value = '  +8.06%  '
$ws.Range('E21').Value = '  +2.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.25'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('E23').Value = '  +2.29%  '
$ws.Range('E24').Value = '  +10.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.10'
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('E27').Value = '  -5.35%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('E30').Value = '  +1.91%  '
$ws.Range('E31').Value = '  +0.93%  '
$ws.Range('E32').Value = '  +2.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.94'
$ws.Range('E33').Value = '  -0.84%  '
$ws.Range('E34').Value = '  +1.18%  '
$ws.Range('E35').Value = '  +1.24%  '
$ws.Range('E36').Value = '  +11.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.125.68'
$ws.Range('E37').Value = '  +3.82%  '
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('E39').Value = '  +2.97%  '
$ws.Range('E40').Value = '  +1.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.489'
$ws.Range('E41').Value = '  -0.86%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.781'
$ws.Range('E42').Value = '  -2.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.734.92'
$ws.Range('E43').Value = '  +2.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.13'
$ws.Range('E44').Value = '  +1.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.87'
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('E46').Value = '  +3.32%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '53.26'
$ws.Range('E47').Value = '  +1.79%  '
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('E49').Value = '  +1.04%  '
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₇0926'
$ws.Range('E51').Value = '  -16.91%  '
